$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# repull data, push all data, mean calculation
# Update the "dSF" (F) column values for the affected rows to reflect
# the repulled/recalculated data.
$ws.Range("F6").Value = 2
$ws.Range("F8").Value = 4
$ws.Range("F11").Value = 4
$ws.Range("F12").Value = 3
$ws.Range("F17").Value = 0
$ws.Range("F18").Value = 1
$ws.Range("F22").Value = -1
$ws.Range("F23").Value = -1
$ws.Range("F29").Value = 4
$ws.Range("F33").Value = 1
$ws.Range("F35").Value = 4
$ws.Range("F39").Value = -2
$ws.Range("F43").Value = -4
